$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# Change 1: "For developers" paragraph - replace the closing sentence.
# -----------------------------------------------------------------------
$d.Content.Find.Execute(
    "for good, and together tries to create an ecosystem as they are growing.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "for good, which greatly improves the efficiency of system development and management.",
    2) | Out-Null

# -----------------------------------------------------------------------
# Change 2 & 3: "For students" paragraph.
#   - The paragraph mark (pPr/rPr) font switches from Segoe UI to SimSun.
#   - The big run is rewritten/split into several runs with new text and
#     mixed Segoe UI / SimSun fonts.
# -----------------------------------------------------------------------
$p = $d.Paragraphs.Item(6)
$r = $p.Range

$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:spacing w:before="100" w:after="240" w:line="360"/><w:ind w:right="0" w:left="0" w:firstLine="0"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="SimSun" w:hAnsi="SimSun" w:cs="SimSun" w:eastAsia="SimSun"/><w:color w:val="333333"/><w:spacing w:val="0"/><w:position w:val="0"/><w:sz w:val="24"/><w:shd w:fill="FFFFFF" w:val="clear"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI" w:eastAsia="Segoe UI"/><w:b/><w:color w:val="333333"/><w:spacing w:val="0"/><w:position w:val="0"/><w:sz w:val="24"/><w:shd w:fill="FFFFFF" w:val="clear"/></w:rPr><w:t xml:space="preserve">For students</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI" w:eastAsia="Segoe UI"/><w:color w:val="333333"/><w:spacing w:val="0"/><w:position w:val="0"/><w:sz w:val="24"/><w:shd w:fill="FFFFFF" w:val="clear"/></w:rPr><w:t xml:space="preserve">, rDSN provides a platform where you can easily simplify, understand and manipulate a distributed system. When learning distributed protocols, you can easily implement one atop of rDSN, and test it on its simulator. The simulator can abstract away many practical difficulties initially, and you can add them back gradually to evolve your protocol, such as from single-thread to multiple-thread, from constant message delay to variant ones, even with message lost. To understand the running protocol, </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Segoe UI" w:hAnsi="Segoe UI" w:cs="Segoe UI" w:eastAsia="Segoe UI"/><w:color w:val="333333"/><w:spacing w:val="0"/><w:position w:val="0"/><w:sz w:val="24"/><w:shd w:fill="FFFFFF" w:val="clear"/></w:rPr><w:t xml:space="preserve">rDSN provides flow tracing and generates a so-called “event matrix” which records the invocation count among different events, revealing the dependencies with weight inside the system.  </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="SimSun" w:hAnsi="SimSun" w:cs="SimSun" w:eastAsia="SimSun"/><w:color w:val="333333"/><w:spacing w:val="0"/><w:position w:val="0"/><w:sz w:val="24"/><w:shd w:fill="FFFFFF" w:val="clear"/></w:rPr><w:t xml:space="preserve">Further</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="SimSun" w:hAnsi="SimSun" w:cs="SimSun" w:eastAsia="SimSun"/><w:color w:val="333333"/><w:spacing w:val="0"/><w:position w:val="0"/><w:sz w:val="24"/><w:shd w:fill="FFFFFF" w:val="clear"/></w:rPr><w:t xml:space="preserve">more, you can easily replace a low component to see what happens. </w:t></w:r></w:p>
'@

$r.InsertXML($xml)

# InsertXML collapses an all-zero <w:ind> element; force it back explicitly
# so the paragraph-format indentation stays present in the saved XML.
$p2 = $d.Paragraphs.Item(6)
$p2.LeftIndent = 0
$p2.RightIndent = 0
$p2.FirstLineIndent = 0
